# "more adjustments to ACH"
# Append 7 new check-register rows (370-376) for the 6/13/2024 ACH batch,
# matching the existing "Check Register" sheet's layout:
#   A = Check #, B = Date, C = Payee, D = Cash Account, E = Amount

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Check Register")

$rows = @(
    @{ Check = "12922"; Date = 45456; Payee = "Bhate Environmental Associates, Inc"; Account = "11040"; Amount = 18043.75 },
    @{ Check = "12923"; Date = 45456; Payee = "Charles County Government";          Account = "11040"; Amount = 91.86 },
    @{ Check = "12924"; Date = 45456; Payee = "Employment Screening Services, Inc"; Account = "11040"; Amount = 67.03 },
    @{ Check = "12925"; Date = 45456; Payee = "Labelmaster";                        Account = "11040"; Amount = 541.98 },
    @{ Check = "12926"; Date = 45456; Payee = "Office Equipment Svcs";              Account = "11040"; Amount = 668.22 },
    @{ Check = "12927"; Date = 45456; Payee = "Safeware, Inc.";                     Account = "11040"; Amount = 4882.91 },
    @{ Check = "12928"; Date = 45456; Payee = "WB Waste-Goode Companies, Inc";      Account = "11040"; Amount = 171 }
)

$startRow = 370
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data.Check
    $ws.Cells.Item($r, 2).Value = $data.Date
    $ws.Cells.Item($r, 3).Value = $data.Payee
    $ws.Cells.Item($r, 4).Value = $data.Account
    $ws.Cells.Item($r, 5).Value = $data.Amount
}

# Keep the frozen header in place, scroll the view down to the new rows,
# and leave the selection where data entry ended.
$ws.Range("A344").Select()
$ws.Range("G371").Select()

Write-Host "Added rows 370-376 to Check Register"
